# run_checker_framework.bat was moved to a new folder; the checker-framework
# results were regenerated, which:
#   - adds a new diagnostic-kind column "[dep-ann] ..." (inserted before the
#     existing "[format.string] ..." column, i.e. at sheet column S)
#   - renumbers three "temp-var-N" diagnostic labels (106->145, 109->148,
#     286->325) in the "[required.method.not.called] ..." columns
#   - adds one new data row "COG Dataset 3 - 7" (inserted right after the
#     existing "COG Dataset 3 - 63" row), with a single count of 1 under the
#     new dep-ann column
#
# We shift cells with plain value copies (bottom/right-to-left) instead of
# Range.Insert/EntireColumn.Insert/EntireRow.Insert so that no incidental
# extra cell-style gets synthesized by a "copy formatting from neighbor"
# side effect - every cell in the grid already carries the style it needs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Phase 1: make room for the new column. Shift columns S(19)..Z(26) one
# slot right, to T(20)..AA(27), across every existing row (1..24).
# Walk right-to-left so we never clobber a source cell before reading it.
# ---------------------------------------------------------------------
for ($col = 26; $col -ge 19; $col--) {
    for ($row = 1; $row -le 24; $row++) {
        $srcCell = $ws.Cells.Item($row, $col)
        $dstCell = $ws.Cells.Item($row, $col + 1)
        $dstCell.Value = $srcCell.Value2
    }
}

# Clear out column S (19) - it will be repopulated as the new column below.
for ($row = 1; $row -le 24; $row++) {
    $ws.Cells.Item($row, 19).Value = $null
}

# New column header (row 1, col S/19).
$ws.Cells.Item(1, 19).Value = "[dep-ann] deprecated item is not annotated with @Deprecated"

# Rename the three shifted temp-var headers (content rename, same cells
# that the loop above already relocated to columns X/Y/Z = 24/25/26).
$ws.Cells.Item(1, 24).Value = "[required.method.not.called] @MustCall method close may not have been invoked on temp-var-145 or any of its aliases."
$ws.Cells.Item(1, 25).Value = "[required.method.not.called] @MustCall method close may not have been invoked on temp-var-148 or any of its aliases."
$ws.Cells.Item(1, 26).Value = "[required.method.not.called] @MustCall method close may not have been invoked on temp-var-325 or any of its aliases."

# ---------------------------------------------------------------------
# Phase 2: make room for the new row. Shift rows 17..24 one slot down,
# to 18..25, across every column (A(1)..AA(27)) in the now-widened grid.
# Walk bottom-to-top so we never clobber a source cell before reading it.
# ---------------------------------------------------------------------
for ($row = 24; $row -ge 17; $row--) {
    for ($col = 1; $col -le 27; $col++) {
        $srcCell = $ws.Cells.Item($row, $col)
        $dstCell = $ws.Cells.Item($row + 1, $col)
        $dstCell.Value = $srcCell.Value2
    }
}

# Clear out row 17 - it will be repopulated as the new row below.
for ($col = 1; $col -le 27; $col++) {
    $ws.Cells.Item(17, $col).Value = $null
}

# New row label (col A) and its single data point under the new dep-ann
# column (col S/19).
$ws.Cells.Item(17, 1).Value = "COG Dataset 3 - 7"
$ws.Cells.Item(17, 19).Value = 1

# Row-1 header style ("A1" style, bold+border+center) needs to be carried
# onto the new A17 label cell, matching every other label cell in col A.
$ws.Cells.Item(16, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The newly-created trailing column (AA, header row only) and the
# newly-created trailing row (25, label cell only) are brand new cells
# that never inherited the shared header/label style via the Value2
# copies above - stamp the same style onto them explicitly.
$ws.Cells.Item(1, 26).Copy()
$ws.Cells.Item(1, 27).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(24, 1).Copy()
$ws.Cells.Item(25, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
